$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.378.67'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.871.33'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.18%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.09%  '
$ws.Range('E7').Value = '  -1.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2868'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06488'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.99%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.78'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '100.56'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.86%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07795'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.872.32'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7246'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.168'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.80%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '284.34'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.368.76'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.06'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.0000'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007479'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.116.02'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.332'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.310'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.99%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.030'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.94'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.892'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.09655'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.45%  '
$ws.Range('E30').Value = '  -1.76%  '
$ws.Range('E31').Value = '  -1.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.225'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.141'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04806'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.123'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6874'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.723'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01897'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.834'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '76.18'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.277'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.74%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.959'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4214'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8237'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '100.76'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.742'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.001'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.98'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.52%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05763'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '886.93'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.89%  '
